$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Fecha (D), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) for rows 2, 3 and 5.
# The changes reflect a row re-ordering (by date) of the weekly price data.

# Row 2
$ws.Range("D2").Value = 44980
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 25000
$ws.Range("O2").Value = 25000
$ws.Range("P2").Value = 25000
$ws.Range("S2").Value = 3125

# Row 3
$ws.Range("D3").Value = 44971
$ws.Range("M3").Value = 25

# Row 5
$ws.Range("D5").Value = 44973
$ws.Range("M5").Value = 55
$ws.Range("N5").Value = 28000
$ws.Range("O5").Value = 28000
$ws.Range("P5").Value = 28000
$ws.Range("S5").Value = 3500
